# maj ipt emilien TP13 et cours 4-3
# Adds 3 new exercise rows (sql / pont_de_wheastone / polynome) right before
# the trailing "fin" marker row of the "exos" sheet, and updates the
# selection/view state accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with a single-column "fin" row at row 113.
# Insert 3 blank rows above it so the new data lands on rows 113-115 and
# "fin" is pushed down to row 116.
$ws.Rows.Item(113).Insert()
$ws.Rows.Item(113).Insert()
$ws.Rows.Item(113).Insert()

# Set the "sql" entry first so its strings are appended to the shared
# string table before the other two new rows (matches original authoring
# order), even though it ends up visually below them on row 115.
$ws.Range("A115").Value = "sql"
$ws.Range("B115").Value = "SQL-004"
$ws.Range("C115").Value = "Base de données sur les acteurs"

$ws.Range("A113").Value = "systemes"
$ws.Range("B113").Value = "pont_de_wheastone"
$ws.Range("C113").Value = "Application de physique : pont de wheastone"

$ws.Range("A114").Value = "systemes"
$ws.Range("B114").Value = "polynome"
$ws.Range("C114").Value = "Interpolation"

# Update the view: scroll so row 97 is at the top, and select C114 (the
# "Interpolation" title cell), matching the authored selection state.
[void]$excel.Goto($ws.Range("A97"), $true)
[void]$ws.Range("C114").Select()
